# Matriz de trazabilidad - actualizacion
# - Actualiza fechas de estado (I) de 08/11/2020 (44143) a 11/11/2020 (44146)
#   en las filas 8-15, 30, 35-39.
# - Marca la fila 10 (columna H "Estado") como "Completado" (antes "Falta"),
#   clonando el formato (relleno/fuente verdes) usado en el resto de filas
#   "Completado" (p.ej. H16).
# - Ajusta el alto de la fila 8.
# - Actualiza la celda activa / seleccion de la hoja a I11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Posicion/tamano de la ventana del libro (best effort; algunos hosts
#     headless no serializan esto en bookViews, pero se fija igualmente) ---
$win = $excel.ActiveWindow
$win.Left = 4545
$win.Top = 1380
$win.Width = 21600
$win.Height = 13470

# --- Fila 8: alto de fila ---
$ws.Rows.Item(8).RowHeight = 97.5

# --- Columna I: fechas de estado 44143 -> 44146 ---
$dateRows = @(8, 9, 10, 11, 12, 13, 14, 15, 30, 35, 36, 37, 38, 39)
foreach ($r in $dateRows) {
    $ws.Cells.Item($r, 9).Value = 44146
}

# --- H10: Estado "Falta" -> "Completado", con el formato verde usado en
#     las demas celdas "Completado" (p.ej. H16) ---
$ws.Range("H16").Copy()
$ws.Range("H10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H10").Value = "Completado"

# --- Desplazamiento de la vista (best effort, topLeftCell "A13") y
#     seleccion activa: I39 -> I11 ---
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("I11").Select() | Out-Null
